# irrigation_main.xlsx -- add the "irrigation_requirement" analysis sheet
# and shift the active-tab / selection bookkeeping accordingly.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Move the selection on "Crop Calendars" before we lose focus on it
#    (matches the diff: selection activeCell moves from E12 to C21 and
#    the tab is no longer the active one once the new sheet is added).
# ------------------------------------------------------------------
$wsCrop = $wb.Worksheets.Item("Crop Calendars")
$wsCrop.Range("C21").Select()

# ------------------------------------------------------------------
# 2. Add the new worksheet as the last tab, named "irrigation_requirement".
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "irrigation_requirement"

# Column widths
$ws.Columns.Item(1).ColumnWidth = 24
$ws.Columns.Item(2).ColumnWidth = 19.5
$ws.Columns.Item(3).ColumnWidth = 24.5
$ws.Columns.Item(4).ColumnWidth = 28
$ws.Columns.Item(5).ColumnWidth = 29
$ws.Columns.Item(6).ColumnWidth = 18.83

# ------------------------------------------------------------------
# 3. Header / notes block
# ------------------------------------------------------------------
$ws.Range("A1").Value = "Source: Phillips Report Appendix A page 2"

$ws.Range("A2").Value = '"Kolars and Mitchell states that after taking into acount potential evapotranspiration with losses that is amount withdrawn (2.5 times PE) and the return flow (35% of the amount withdrawn) the irrigation requirement becomes approximately 1.6 m (1*2.5 -(2.5*0.35)=1.6). Asuming the same distribution for the Tigris Basin with corrected consumptive use:" irrigation requirement unit is meters/month'
$ws.Range("A2").Font.Name = "ArialMT"
$ws.Range("A2").Font.Size = 10

$ws.Range("A3").Font.Name = "ArialMT"
$ws.Range("A3").Font.Size = 10

# ------------------------------------------------------------------
# 4. Area / conversion block (rows 4-9)
# ------------------------------------------------------------------
$ws.Range("B4").Value = "irrigated_area_ha"
$ws.Range("C4").Value = "irrigated_area_m2"
$ws.Range("E4").Value = "ha2squaremeters"
$ws.Range("F4").Value = 10000

$ws.Range("A5").Value = "diversion_downstream_planned"
$ws.Range("B5").Value = 121000
$ws.Range("C5").Formula = "=B5*`$F`$4"
$ws.Range("E5").Value = "seconds_per_month"
$ws.Range("F5").Formula = "=60*60*24*31"

$ws.Range("A6").Value = "diversion_upstream_existing"
$ws.Range("B6").Value = 138000
$ws.Range("C6").Formula = "=B6*`$F`$4"

$ws.Range("A7").Value = "diversion_upstream_planned"
$ws.Range("B7").Value = 375000
$ws.Range("C7").Formula = "=B7*`$F`$4"

$ws.Range("A8").Value = "diversion_upstream_total"
$ws.Range("B8").Formula = "=SUM(B6:B7)"
$ws.Range("C8").Formula = "=SUM(C6:C7)"

$ws.Range("A9").Value = "diversion_sum"
$ws.Range("B9").Formula = "=SUM(B5,B8)"
$ws.Range("C9").Formula = "=SUM(C5,C8)"
$ws.Range("D9").Value = "of all the planned and existing irrigation projects up- and downstream of Ilisu dam"

# ------------------------------------------------------------------
# 5. Units + monthly table headers (rows 13-14)
# ------------------------------------------------------------------
$ws.Range("A13").Value = "units"
$ws.Range("B13").Value = "m/month"
$ws.Range("C13").Value = "m^3/s"
$ws.Range("D13").Value = "m^3/s"
$ws.Range("E13").Value = "m^3/s"
$ws.Range("F13").Value = "m^3/s"

$ws.Range("A14").Value = "month"
$ws.Range("B14").Value = "irrigation_requirement"
$ws.Range("C14").Value = "diversion_upstream_existing"
$ws.Range("D14").Value = "diversion_upstream_planned"
$ws.Range("E14").Value = "diversion_downstream_planned"
$ws.Range("F14").Value = "diversion_sum"

# ------------------------------------------------------------------
# 6. Monthly data table (rows 15-26)
# ------------------------------------------------------------------
$months = 1..12
$irrigReq = @(0, 0, 0, 0.066000000000000003, 0.13300000000000001, 0.33400000000000002, `
              0.46200000000000002, 0.39, 0.187, 0.027, 0, 0)

for ($i = 0; $i -lt 12; $i++) {
    $r = 15 + $i
    $ws.Cells.Item($r, 1).Value = $months[$i]
    $bCell = $ws.Cells.Item($r, 2)
    $bCell.Value = $irrigReq[$i]
    $bCell.NumberFormat = "0.000"

    $cCell = $ws.Cells.Item($r, 3)
    $cCell.Formula = "=B$r*`$C`$6/`$F`$5"
    $cCell.NumberFormat = "0.000"

    $dCell = $ws.Cells.Item($r, 4)
    $dCell.Formula = "=B$r*`$C`$7/`$F`$5"
    $dCell.NumberFormat = "0.000"

    $eCell = $ws.Cells.Item($r, 5)
    $eCell.Formula = "=B$r*`$C`$5/`$F`$5"
    $eCell.NumberFormat = "0.000"

    $fCell = $ws.Cells.Item($r, 6)
    $fCell.Formula = "=SUM(C$r`:E$r)"
    $fCell.NumberFormat = "0.000"
}

# ------------------------------------------------------------------
# 7. Leave the selection on the new sheet matching the diff (B12) and
#    make sure it is the active sheet / tab.
# ------------------------------------------------------------------
$ws.Range("B12").Select()

Write-Output "irrigation_requirement sheet created"
